$d = $word.ActiveDocument

# --- Title date line -------------------------------------------------
$d.Content.Find.Execute("2024-11-27 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-11-28 Thursday", 2)

# --- Table of multiplication problems --------------------------------
$t = $d.Tables.Item(1)

$cellMap = @{
    1  = @{1="85×73=6205"; 2="44×25=1100"; 3="98×71=6958"; 4="16×91=1456"; 5="69×65=4485"}
    5  = @{1="82×68=5576"; 2="30×71=2130"; 3="47×49=2303"; 4="87×14=1218"; 5="67×53=3551"}
    10 = @{1="31×43=1333"; 2="90×61=5490"; 3="24×13=312";  4="19×37=703";  5="62×74=4588"}
    15 = @{1="58×19=1102"; 2="58×67=3886"; 3="12×59=708";  4="58×27=1566"; 5="16×37=592"}
    20 = @{1="90×30=2700"; 2="71×70=4970"; 3="81×79=6399"; 4="23×92=2116"; 5="55×60=3300"}
}

foreach ($rowIndex in $cellMap.Keys) {
    $row = $cellMap[$rowIndex]
    foreach ($colIndex in $row.Keys) {
        $newText = $row[$colIndex]
        $cell = $t.Cell($rowIndex, $colIndex)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $newText
    }
}
